# Edit: B1--and-B2-PowerPoint.pptx
#  1) Slide 5's table switches from the custom "Table_0" table style to the
#     built-in "No Style, No Grid" table style.
#  2) The deck's applied theme ("Integral" / Red Violet) is swapped back to
#     the default "Office Theme" colors.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 (the table is the 2nd shape on the slide) ---
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{DB244FB5-A3C7-4ADC-B93B-C65F3AE9C532}")

# --- 2) Recolor the active theme to the standard "Office" palette ---
function HexToColorRef($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Order matches ThemeColorScheme.Item(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToColorRef $officeColors[$i - 1]
}
